$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 2010年
$ws.Cells.Item(2, 1).Value = "2010年"
$ws.Cells.Item(2, 2).Value = 11
$ws.Cells.Item(2, 3).Value = 15
$ws.Cells.Item(2, 4).Value = 8
$ws.Cells.Item(2, 5).Value = 6
$ws.Cells.Item(2, 6).Value = 4
$ws.Cells.Item(2, 7).Value = 14
$ws.Cells.Item(2, 8).Value = 15
$ws.Cells.Item(2, 9).Value = 12
$ws.Cells.Item(2, 10).Value = 15
$ws.Cells.Item(2, 11).Value = 8
$ws.Cells.Item(2, 12).Value = 6
$ws.Cells.Item(2, 13).Value = 4
$ws.Cells.Item(2, 14).Value = 18
$ws.Cells.Item(2, 15).Value = 16
$ws.Cells.Item(2, 16).Value = 19
$ws.Cells.Item(2, 17).Value = 20
$ws.Cells.Item(2, 18).Value = 1
$ws.Cells.Item(2, 19).Value = 2
$ws.Cells.Item(2, 20).Value = 5
$ws.Cells.Item(2, 21).Value = 9
$ws.Cells.Item(2, 22).Value = 7
$ws.Cells.Item(2, 23).Value = 10
$ws.Cells.Item(2, 24).Value = 13
$ws.Cells.Item(2, 25).Value = 14
$ws.Cells.Item(2, 26).Value = 3
$ws.Cells.Item(2, 27).Value = 17
$ws.Cells.Item(2, 28).Value = 11
$ws.Cells.Item(2, 29).Value = 18
$ws.Cells.Item(2, 30).Value = 18
$ws.Cells.Item(2, 31).Value = 2
$ws.Cells.Item(2, 32).Value = 1
$ws.Cells.Item(2, 33).Value = 5
$ws.Cells.Item(2, 34).Value = 9
$ws.Cells.Item(2, 35).Value = 7
$ws.Cells.Item(2, 36).Value = 11
$ws.Cells.Item(2, 37).Value = 16
$ws.Cells.Item(2, 38).Value = 9
$ws.Cells.Item(2, 39).Value = 7
$ws.Cells.Item(2, 40).Value = 4
$ws.Cells.Item(2, 41).Value = 13
$ws.Cells.Item(2, 42).Value = 15
$ws.Cells.Item(2, 43).Value = 19
$ws.Cells.Item(2, 44).Value = 3
$ws.Cells.Item(2, 45).Value = 1
$ws.Cells.Item(2, 46).Value = 5
$ws.Cells.Item(2, 47).Value = 8
$ws.Cells.Item(2, 48).Value = 6
$ws.Cells.Item(2, 49).Value = 10
$ws.Cells.Item(2, 50).Value = 14
$ws.Cells.Item(2, 51).Value = 18
$ws.Cells.Item(2, 52).Value = 2
$ws.Cells.Item(2, 53).Value = 17
$ws.Cells.Item(2, 54).Value = 12
$ws.Cells.Item(2, 55).Value = 10
$ws.Cells.Item(2, 56).Value = 13
$ws.Cells.Item(2, 57).Value = 16
$ws.Cells.Item(2, 58).Value = 3
$ws.Cells.Item(2, 59).Value = 17
$ws.Cells.Item(2, 60).Value = 12

# Row 3: 2011年
$ws.Cells.Item(3, 1).Value = "2011年"
$ws.Cells.Item(3, 2).Value = 11
$ws.Cells.Item(3, 3).Value = 16
$ws.Cells.Item(3, 4).Value = 8
$ws.Cells.Item(3, 5).Value = 6
$ws.Cells.Item(3, 6).Value = 4
$ws.Cells.Item(3, 7).Value = 14
$ws.Cells.Item(3, 8).Value = 15
$ws.Cells.Item(3, 9).Value = 13
$ws.Cells.Item(3, 10).Value = 15
$ws.Cells.Item(3, 11).Value = 8
$ws.Cells.Item(3, 12).Value = 6
$ws.Cells.Item(3, 13).Value = 4
$ws.Cells.Item(3, 14).Value = 18
$ws.Cells.Item(3, 15).Value = 16
$ws.Cells.Item(3, 16).Value = 20
$ws.Cells.Item(3, 17).Value = 19
$ws.Cells.Item(3, 18).Value = 2
$ws.Cells.Item(3, 19).Value = 1
$ws.Cells.Item(3, 20).Value = 5
$ws.Cells.Item(3, 21).Value = 10
$ws.Cells.Item(3, 22).Value = 7
$ws.Cells.Item(3, 23).Value = 9
$ws.Cells.Item(3, 24).Value = 11
$ws.Cells.Item(3, 25).Value = 14
$ws.Cells.Item(3, 26).Value = 3
$ws.Cells.Item(3, 27).Value = 17
$ws.Cells.Item(3, 28).Value = 12
$ws.Cells.Item(3, 29).Value = 20
$ws.Cells.Item(3, 30).Value = 19
$ws.Cells.Item(3, 31).Value = 2
$ws.Cells.Item(3, 32).Value = 1
$ws.Cells.Item(3, 33).Value = 5
$ws.Cells.Item(3, 34).Value = 10
$ws.Cells.Item(3, 35).Value = 7
$ws.Cells.Item(3, 36).Value = 11
$ws.Cells.Item(3, 37).Value = 16
$ws.Cells.Item(3, 38).Value = 9
$ws.Cells.Item(3, 39).Value = 7
$ws.Cells.Item(3, 40).Value = 4
$ws.Cells.Item(3, 41).Value = 13
$ws.Cells.Item(3, 42).Value = 15
$ws.Cells.Item(3, 43).Value = 19
$ws.Cells.Item(3, 44).Value = 2
$ws.Cells.Item(3, 45).Value = 1
$ws.Cells.Item(3, 46).Value = 5
$ws.Cells.Item(3, 47).Value = 10
$ws.Cells.Item(3, 48).Value = 6
$ws.Cells.Item(3, 49).Value = 8
$ws.Cells.Item(3, 50).Value = 14
$ws.Cells.Item(3, 51).Value = 18
$ws.Cells.Item(3, 52).Value = 3
$ws.Cells.Item(3, 53).Value = 17
$ws.Cells.Item(3, 54).Value = 12
$ws.Cells.Item(3, 55).Value = 9
$ws.Cells.Item(3, 56).Value = 13
$ws.Cells.Item(3, 57).Value = 17
$ws.Cells.Item(3, 58).Value = 3
$ws.Cells.Item(3, 59).Value = 18
$ws.Cells.Item(3, 60).Value = 12

# Row 4: 2012年
$ws.Cells.Item(4, 1).Value = "2012年"
$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 3).Value = 16
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(4, 7).Value = 14
$ws.Cells.Item(4, 8).Value = 15
$ws.Cells.Item(4, 9).Value = 11
$ws.Cells.Item(4, 10).Value = 15
$ws.Cells.Item(4, 11).Value = 5
$ws.Cells.Item(4, 12).Value = 7
$ws.Cells.Item(4, 13).Value = 4
$ws.Cells.Item(4, 14).Value = 18
$ws.Cells.Item(4, 15).Value = 16
$ws.Cells.Item(4, 16).Value = 19
$ws.Cells.Item(4, 17).Value = 20
$ws.Cells.Item(4, 18).Value = 1
$ws.Cells.Item(4, 19).Value = 2
$ws.Cells.Item(4, 20).Value = 6
$ws.Cells.Item(4, 21).Value = 10
$ws.Cells.Item(4, 22).Value = 8
$ws.Cells.Item(4, 23).Value = 9
$ws.Cells.Item(4, 24).Value = 12
$ws.Cells.Item(4, 25).Value = 13
$ws.Cells.Item(4, 26).Value = 3
$ws.Cells.Item(4, 27).Value = 17
$ws.Cells.Item(4, 28).Value = 14
$ws.Cells.Item(4, 29).Value = 20
$ws.Cells.Item(4, 30).Value = 19
$ws.Cells.Item(4, 31).Value = 2
$ws.Cells.Item(4, 32).Value = 1
$ws.Cells.Item(4, 33).Value = 5
$ws.Cells.Item(4, 34).Value = 10
$ws.Cells.Item(4, 35).Value = 8
$ws.Cells.Item(4, 36).Value = 11
$ws.Cells.Item(4, 37).Value = 15
$ws.Cells.Item(4, 38).Value = 6
$ws.Cells.Item(4, 39).Value = 8
$ws.Cells.Item(4, 40).Value = 4
$ws.Cells.Item(4, 41).Value = 13
$ws.Cells.Item(4, 42).Value = 14
$ws.Cells.Item(4, 43).Value = 19
$ws.Cells.Item(4, 44).Value = 2
$ws.Cells.Item(4, 45).Value = 1
$ws.Cells.Item(4, 46).Value = 5
$ws.Cells.Item(4, 47).Value = 10
$ws.Cells.Item(4, 48).Value = 7
$ws.Cells.Item(4, 49).Value = 9
$ws.Cells.Item(4, 50).Value = 16
$ws.Cells.Item(4, 51).Value = 18
$ws.Cells.Item(4, 52).Value = 3
$ws.Cells.Item(4, 53).Value = 17
$ws.Cells.Item(4, 54).Value = 12
$ws.Cells.Item(4, 55).Value = 9
$ws.Cells.Item(4, 56).Value = 13
$ws.Cells.Item(4, 57).Value = 17
$ws.Cells.Item(4, 58).Value = 3
$ws.Cells.Item(4, 59).Value = 18
$ws.Cells.Item(4, 60).Value = 12

# Row 5: 2013年
$ws.Cells.Item(5, 1).Value = "2013年"
$ws.Cells.Item(5, 2).Value = 8
$ws.Cells.Item(5, 3).Value = 13
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = 6
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = ""
$ws.Cells.Item(5, 8).Value = 12
$ws.Cells.Item(5, 9).Value = 10
$ws.Cells.Item(5, 10).Value = 13
$ws.Cells.Item(5, 11).Value = ""
$ws.Cells.Item(5, 12).Value = 6
$ws.Cells.Item(5, 13).Value = 4
$ws.Cells.Item(5, 14).Value = ""
$ws.Cells.Item(5, 15).Value = 14
$ws.Cells.Item(5, 16).Value = 16
$ws.Cells.Item(5, 17).Value = 17
$ws.Cells.Item(5, 18).Value = 1
$ws.Cells.Item(5, 19).Value = 2
$ws.Cells.Item(5, 20).Value = 5
$ws.Cells.Item(5, 21).Value = 9
$ws.Cells.Item(5, 22).Value = 7
$ws.Cells.Item(5, 23).Value = 8
$ws.Cells.Item(5, 24).Value = 11
$ws.Cells.Item(5, 25).Value = 12
$ws.Cells.Item(5, 26).Value = 3
$ws.Cells.Item(5, 27).Value = 15
$ws.Cells.Item(5, 28).Value = ""
$ws.Cells.Item(5, 29).Value = 16
$ws.Cells.Item(5, 30).Value = 17
$ws.Cells.Item(5, 31).Value = 2
$ws.Cells.Item(5, 32).Value = 1
$ws.Cells.Item(5, 33).Value = 5
$ws.Cells.Item(5, 34).Value = 10
$ws.Cells.Item(5, 35).Value = 7
$ws.Cells.Item(5, 36).Value = 8
$ws.Cells.Item(5, 37).Value = 13
$ws.Cells.Item(5, 38).Value = ""
$ws.Cells.Item(5, 39).Value = 7
$ws.Cells.Item(5, 40).Value = 4
$ws.Cells.Item(5, 41).Value = ""
$ws.Cells.Item(5, 42).Value = 12
$ws.Cells.Item(5, 43).Value = 16
$ws.Cells.Item(5, 44).Value = 3
$ws.Cells.Item(5, 45).Value = 1
$ws.Cells.Item(5, 46).Value = 5
$ws.Cells.Item(5, 47).Value = 10
$ws.Cells.Item(5, 48).Value = 6
$ws.Cells.Item(5, 49).Value = 9
$ws.Cells.Item(5, 50).Value = 11
$ws.Cells.Item(5, 51).Value = 15
$ws.Cells.Item(5, 52).Value = 2
$ws.Cells.Item(5, 53).Value = 14
$ws.Cells.Item(5, 54).Value = ""
$ws.Cells.Item(5, 55).Value = 9
$ws.Cells.Item(5, 56).Value = 11
$ws.Cells.Item(5, 57).Value = 14
$ws.Cells.Item(5, 58).Value = 3
$ws.Cells.Item(5, 59).Value = 15
$ws.Cells.Item(5, 60).Value = ""

# Row 6: 2014年
$ws.Cells.Item(6, 1).Value = "2014年"
$ws.Cells.Item(6, 2).Value = 10
$ws.Cells.Item(6, 3).Value = 13
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 5).Value = 6
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = 12
$ws.Cells.Item(6, 9).Value = 10
$ws.Cells.Item(6, 10).Value = 14
$ws.Cells.Item(6, 11).Value = ""
$ws.Cells.Item(6, 12).Value = 6
$ws.Cells.Item(6, 13).Value = 4
$ws.Cells.Item(6, 14).Value = ""
$ws.Cells.Item(6, 15).Value = 13
$ws.Cells.Item(6, 16).Value = 16
$ws.Cells.Item(6, 17).Value = 17
$ws.Cells.Item(6, 18).Value = 1
$ws.Cells.Item(6, 19).Value = 2
$ws.Cells.Item(6, 20).Value = 5
$ws.Cells.Item(6, 21).Value = 9
$ws.Cells.Item(6, 22).Value = 7
$ws.Cells.Item(6, 23).Value = 8
$ws.Cells.Item(6, 24).Value = 11
$ws.Cells.Item(6, 25).Value = 12
$ws.Cells.Item(6, 26).Value = 3
$ws.Cells.Item(6, 27).Value = 15
$ws.Cells.Item(6, 28).Value = ""
$ws.Cells.Item(6, 29).Value = 16
$ws.Cells.Item(6, 30).Value = 17
$ws.Cells.Item(6, 31).Value = 2
$ws.Cells.Item(6, 32).Value = 1
$ws.Cells.Item(6, 33).Value = 5
$ws.Cells.Item(6, 34).Value = 9
$ws.Cells.Item(6, 35).Value = 4
$ws.Cells.Item(6, 36).Value = 8
$ws.Cells.Item(6, 37).Value = 13
$ws.Cells.Item(6, 38).Value = ""
$ws.Cells.Item(6, 39).Value = 7
$ws.Cells.Item(6, 40).Value = 4
$ws.Cells.Item(6, 41).Value = ""
$ws.Cells.Item(6, 42).Value = 12
$ws.Cells.Item(6, 43).Value = 16
$ws.Cells.Item(6, 44).Value = 2
$ws.Cells.Item(6, 45).Value = 1
$ws.Cells.Item(6, 46).Value = 5
$ws.Cells.Item(6, 47).Value = 9
$ws.Cells.Item(6, 48).Value = 6
$ws.Cells.Item(6, 49).Value = 10
$ws.Cells.Item(6, 50).Value = 11
$ws.Cells.Item(6, 51).Value = 15
$ws.Cells.Item(6, 52).Value = 3
$ws.Cells.Item(6, 53).Value = 14
$ws.Cells.Item(6, 54).Value = ""
$ws.Cells.Item(6, 55).Value = 8
$ws.Cells.Item(6, 56).Value = 11
$ws.Cells.Item(6, 57).Value = 14
$ws.Cells.Item(6, 58).Value = 3
$ws.Cells.Item(6, 59).Value = 15
$ws.Cells.Item(6, 60).Value = ""

# Row 7: 2015年
$ws.Cells.Item(7, 1).Value = "2015年"
$ws.Cells.Item(7, 2).Value = 9
$ws.Cells.Item(7, 3).Value = 13
$ws.Cells.Item(7, 4).Value = ""
$ws.Cells.Item(7, 5).Value = 6
$ws.Cells.Item(7, 6).Value = 4
$ws.Cells.Item(7, 7).Value = ""
$ws.Cells.Item(7, 8).Value = 14
$ws.Cells.Item(7, 9).Value = 10
$ws.Cells.Item(7, 10).Value = 13
$ws.Cells.Item(7, 11).Value = ""
$ws.Cells.Item(7, 12).Value = 6
$ws.Cells.Item(7, 13).Value = 4
$ws.Cells.Item(7, 14).Value = ""
$ws.Cells.Item(7, 15).Value = 14
$ws.Cells.Item(7, 16).Value = 16
$ws.Cells.Item(7, 17).Value = 17
$ws.Cells.Item(7, 18).Value = 1
$ws.Cells.Item(7, 19).Value = 2
$ws.Cells.Item(7, 20).Value = 5
$ws.Cells.Item(7, 21).Value = 9
$ws.Cells.Item(7, 22).Value = 7
$ws.Cells.Item(7, 23).Value = 8
$ws.Cells.Item(7, 24).Value = 11
$ws.Cells.Item(7, 25).Value = 12
$ws.Cells.Item(7, 26).Value = 3
$ws.Cells.Item(7, 27).Value = 15
$ws.Cells.Item(7, 28).Value = ""
$ws.Cells.Item(7, 29).Value = 16
$ws.Cells.Item(7, 30).Value = 17
$ws.Cells.Item(7, 31).Value = 2
$ws.Cells.Item(7, 32).Value = 1
$ws.Cells.Item(7, 33).Value = 5
$ws.Cells.Item(7, 34).Value = 10
$ws.Cells.Item(7, 35).Value = 7
$ws.Cells.Item(7, 36).Value = 8
$ws.Cells.Item(7, 37).Value = 13
$ws.Cells.Item(7, 38).Value = ""
$ws.Cells.Item(7, 39).Value = 6
$ws.Cells.Item(7, 40).Value = 4
$ws.Cells.Item(7, 41).Value = ""
$ws.Cells.Item(7, 42).Value = 12
$ws.Cells.Item(7, 43).Value = 16
$ws.Cells.Item(7, 44).Value = 3
$ws.Cells.Item(7, 45).Value = 1
$ws.Cells.Item(7, 46).Value = 5
$ws.Cells.Item(7, 47).Value = 9
$ws.Cells.Item(7, 48).Value = 7
$ws.Cells.Item(7, 49).Value = 10
$ws.Cells.Item(7, 50).Value = 11
$ws.Cells.Item(7, 51).Value = 14
$ws.Cells.Item(7, 52).Value = 2
$ws.Cells.Item(7, 53).Value = 15
$ws.Cells.Item(7, 54).Value = ""
$ws.Cells.Item(7, 55).Value = 8
$ws.Cells.Item(7, 56).Value = 11
$ws.Cells.Item(7, 57).Value = 12
$ws.Cells.Item(7, 58).Value = 3
$ws.Cells.Item(7, 59).Value = 15
$ws.Cells.Item(7, 60).Value = ""

# Row 8: 2016年
$ws.Cells.Item(8, 1).Value = "2016年"
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = 6
$ws.Cells.Item(8, 6).Value = 4
$ws.Cells.Item(8, 7).Value = ""
$ws.Cells.Item(8, 8).Value = 13
$ws.Cells.Item(8, 9).Value = ""
$ws.Cells.Item(8, 10).Value = 14
$ws.Cells.Item(8, 11).Value = ""
$ws.Cells.Item(8, 12).Value = 6
$ws.Cells.Item(8, 13).Value = 4
$ws.Cells.Item(8, 14).Value = ""
$ws.Cells.Item(8, 15).Value = 13
$ws.Cells.Item(8, 16).Value = 16
$ws.Cells.Item(8, 17).Value = 17
$ws.Cells.Item(8, 18).Value = 1
$ws.Cells.Item(8, 19).Value = 2
$ws.Cells.Item(8, 20).Value = 5
$ws.Cells.Item(8, 21).Value = 9
$ws.Cells.Item(8, 22).Value = 7
$ws.Cells.Item(8, 23).Value = 8
$ws.Cells.Item(8, 24).Value = 11
$ws.Cells.Item(8, 25).Value = 12
$ws.Cells.Item(8, 26).Value = 3
$ws.Cells.Item(8, 27).Value = 15
$ws.Cells.Item(8, 28).Value = ""
$ws.Cells.Item(8, 29).Value = 16
$ws.Cells.Item(8, 30).Value = 17
$ws.Cells.Item(8, 31).Value = 2
$ws.Cells.Item(8, 32).Value = 1
$ws.Cells.Item(8, 33).Value = 5
$ws.Cells.Item(8, 34).Value = 9
$ws.Cells.Item(8, 35).Value = 7
$ws.Cells.Item(8, 36).Value = ""
$ws.Cells.Item(8, 37).Value = 14
$ws.Cells.Item(8, 38).Value = ""
$ws.Cells.Item(8, 39).Value = 6
$ws.Cells.Item(8, 40).Value = 4
$ws.Cells.Item(8, 41).Value = ""
$ws.Cells.Item(8, 42).Value = 12
$ws.Cells.Item(8, 43).Value = 16
$ws.Cells.Item(8, 44).Value = 2
$ws.Cells.Item(8, 45).Value = 1
$ws.Cells.Item(8, 46).Value = 5
$ws.Cells.Item(8, 47).Value = 10
$ws.Cells.Item(8, 48).Value = 7
$ws.Cells.Item(8, 49).Value = 9
$ws.Cells.Item(8, 50).Value = 11
$ws.Cells.Item(8, 51).Value = 13
$ws.Cells.Item(8, 52).Value = 3
$ws.Cells.Item(8, 53).Value = 15
$ws.Cells.Item(8, 54).Value = ""
$ws.Cells.Item(8, 55).Value = 8
$ws.Cells.Item(8, 56).Value = 11
$ws.Cells.Item(8, 57).Value = 12
$ws.Cells.Item(8, 58).Value = 3
$ws.Cells.Item(8, 59).Value = 15
$ws.Cells.Item(8, 60).Value = ""

# Row 9: 2017年
$ws.Cells.Item(9, 1).Value = "2017年"
$ws.Cells.Item(9, 2).Value = 10
$ws.Cells.Item(9, 3).Value = 14
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = 6
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(9, 7).Value = ""
$ws.Cells.Item(9, 8).Value = 13
$ws.Cells.Item(9, 9).Value = 10
$ws.Cells.Item(9, 10).Value = 13
$ws.Cells.Item(9, 11).Value = ""
$ws.Cells.Item(9, 12).Value = 6
$ws.Cells.Item(9, 13).Value = 4
$ws.Cells.Item(9, 14).Value = ""
$ws.Cells.Item(9, 15).Value = 14
$ws.Cells.Item(9, 16).Value = 16
$ws.Cells.Item(9, 17).Value = 17
$ws.Cells.Item(9, 18).Value = 1
$ws.Cells.Item(9, 19).Value = 2
$ws.Cells.Item(9, 20).Value = 5
$ws.Cells.Item(9, 21).Value = 9
$ws.Cells.Item(9, 22).Value = 7
$ws.Cells.Item(9, 23).Value = 8
$ws.Cells.Item(9, 24).Value = 12
$ws.Cells.Item(9, 25).Value = 11
$ws.Cells.Item(9, 26).Value = 3
$ws.Cells.Item(9, 27).Value = 15
$ws.Cells.Item(9, 28).Value = ""
$ws.Cells.Item(9, 29).Value = 16
$ws.Cells.Item(9, 30).Value = 17
$ws.Cells.Item(9, 31).Value = 2
$ws.Cells.Item(9, 32).Value = 1
$ws.Cells.Item(9, 33).Value = 5
$ws.Cells.Item(9, 34).Value = 9
$ws.Cells.Item(9, 35).Value = 7
$ws.Cells.Item(9, 36).Value = 8
$ws.Cells.Item(9, 37).Value = 14
$ws.Cells.Item(9, 38).Value = ""
$ws.Cells.Item(9, 39).Value = 6
$ws.Cells.Item(9, 40).Value = 4
$ws.Cells.Item(9, 41).Value = ""
$ws.Cells.Item(9, 42).Value = 12
$ws.Cells.Item(9, 43).Value = 16
$ws.Cells.Item(9, 44).Value = 2
$ws.Cells.Item(9, 45).Value = 1
$ws.Cells.Item(9, 46).Value = 5
$ws.Cells.Item(9, 47).Value = 10
$ws.Cells.Item(9, 48).Value = 7
$ws.Cells.Item(9, 49).Value = 9
$ws.Cells.Item(9, 50).Value = 11
$ws.Cells.Item(9, 51).Value = 13
$ws.Cells.Item(9, 52).Value = 3
$ws.Cells.Item(9, 53).Value = 15
$ws.Cells.Item(9, 54).Value = ""
$ws.Cells.Item(9, 55).Value = 8
$ws.Cells.Item(9, 56).Value = 11
$ws.Cells.Item(9, 57).Value = 12
$ws.Cells.Item(9, 58).Value = 3
$ws.Cells.Item(9, 59).Value = 15
$ws.Cells.Item(9, 60).Value = ""

# Row 10: 2018年
$ws.Cells.Item(10, 1).Value = "2018年"
$ws.Cells.Item(10, 2).Value = 10
$ws.Cells.Item(10, 3).Value = 14
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = 6
$ws.Cells.Item(10, 6).Value = 4
$ws.Cells.Item(10, 7).Value = ""
$ws.Cells.Item(10, 8).Value = 15
$ws.Cells.Item(10, 9).Value = 10
$ws.Cells.Item(10, 10).Value = 14
$ws.Cells.Item(10, 11).Value = ""
$ws.Cells.Item(10, 12).Value = 6
$ws.Cells.Item(10, 13).Value = 4
$ws.Cells.Item(10, 14).Value = ""
$ws.Cells.Item(10, 15).Value = 15
$ws.Cells.Item(10, 16).Value = 16
$ws.Cells.Item(10, 17).Value = 17
$ws.Cells.Item(10, 18).Value = 1
$ws.Cells.Item(10, 19).Value = 2
$ws.Cells.Item(10, 20).Value = 5
$ws.Cells.Item(10, 21).Value = 9
$ws.Cells.Item(10, 22).Value = 7
$ws.Cells.Item(10, 23).Value = 8
$ws.Cells.Item(10, 24).Value = 11
$ws.Cells.Item(10, 25).Value = 12
$ws.Cells.Item(10, 26).Value = 3
$ws.Cells.Item(10, 27).Value = 13
$ws.Cells.Item(10, 28).Value = ""
$ws.Cells.Item(10, 29).Value = 16
$ws.Cells.Item(10, 30).Value = 17
$ws.Cells.Item(10, 31).Value = 2
$ws.Cells.Item(10, 32).Value = 1
$ws.Cells.Item(10, 33).Value = 5
$ws.Cells.Item(10, 34).Value = 9
$ws.Cells.Item(10, 35).Value = 7
$ws.Cells.Item(10, 36).Value = 9
$ws.Cells.Item(10, 37).Value = 14
$ws.Cells.Item(10, 38).Value = ""
$ws.Cells.Item(10, 39).Value = 6
$ws.Cells.Item(10, 40).Value = 4
$ws.Cells.Item(10, 41).Value = ""
$ws.Cells.Item(10, 42).Value = 13
$ws.Cells.Item(10, 43).Value = 16
$ws.Cells.Item(10, 44).Value = 2
$ws.Cells.Item(10, 45).Value = 1
$ws.Cells.Item(10, 46).Value = 5
$ws.Cells.Item(10, 47).Value = 10
$ws.Cells.Item(10, 48).Value = 7
$ws.Cells.Item(10, 49).Value = 8
$ws.Cells.Item(10, 50).Value = 11
$ws.Cells.Item(10, 51).Value = 12
$ws.Cells.Item(10, 52).Value = 3
$ws.Cells.Item(10, 53).Value = 15
$ws.Cells.Item(10, 54).Value = ""
$ws.Cells.Item(10, 55).Value = 8
$ws.Cells.Item(10, 56).Value = 11
$ws.Cells.Item(10, 57).Value = 12
$ws.Cells.Item(10, 58).Value = 3
$ws.Cells.Item(10, 59).Value = 13
$ws.Cells.Item(10, 60).Value = ""

# Row 11: 2019年
$ws.Cells.Item(11, 1).Value = "2019年"
$ws.Cells.Item(11, 2).Value = ""
$ws.Cells.Item(11, 3).Value = 14
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = 6
$ws.Cells.Item(11, 6).Value = 4
$ws.Cells.Item(11, 7).Value = ""
$ws.Cells.Item(11, 8).Value = 15
$ws.Cells.Item(11, 9).Value = ""
$ws.Cells.Item(11, 10).Value = 14
$ws.Cells.Item(11, 11).Value = ""
$ws.Cells.Item(11, 12).Value = 6
$ws.Cells.Item(11, 13).Value = 4
$ws.Cells.Item(11, 14).Value = ""
$ws.Cells.Item(11, 15).Value = 15
$ws.Cells.Item(11, 16).Value = 16
$ws.Cells.Item(11, 17).Value = 17
$ws.Cells.Item(11, 18).Value = 1
$ws.Cells.Item(11, 19).Value = 2
$ws.Cells.Item(11, 20).Value = 5
$ws.Cells.Item(11, 21).Value = 9
$ws.Cells.Item(11, 22).Value = 7
$ws.Cells.Item(11, 23).Value = 8
$ws.Cells.Item(11, 24).Value = 11
$ws.Cells.Item(11, 25).Value = 12
$ws.Cells.Item(11, 26).Value = 3
$ws.Cells.Item(11, 27).Value = 13
$ws.Cells.Item(11, 28).Value = ""
$ws.Cells.Item(11, 29).Value = 17
$ws.Cells.Item(11, 30).Value = 16
$ws.Cells.Item(11, 31).Value = 2
$ws.Cells.Item(11, 32).Value = 1
$ws.Cells.Item(11, 33).Value = 5
$ws.Cells.Item(11, 34).Value = 9
$ws.Cells.Item(11, 35).Value = 7
$ws.Cells.Item(11, 36).Value = ""
$ws.Cells.Item(11, 37).Value = 15
$ws.Cells.Item(11, 38).Value = ""
$ws.Cells.Item(11, 39).Value = 6
$ws.Cells.Item(11, 40).Value = 4
$ws.Cells.Item(11, 41).Value = ""
$ws.Cells.Item(11, 42).Value = 13
$ws.Cells.Item(11, 43).Value = 16
$ws.Cells.Item(11, 44).Value = 2
$ws.Cells.Item(11, 45).Value = 1
$ws.Cells.Item(11, 46).Value = 5
$ws.Cells.Item(11, 47).Value = 10
$ws.Cells.Item(11, 48).Value = 7
$ws.Cells.Item(11, 49).Value = 8
$ws.Cells.Item(11, 50).Value = 11
$ws.Cells.Item(11, 51).Value = 12
$ws.Cells.Item(11, 52).Value = 3
$ws.Cells.Item(11, 53).Value = 13
$ws.Cells.Item(11, 54).Value = ""
$ws.Cells.Item(11, 55).Value = 8
$ws.Cells.Item(11, 56).Value = 11
$ws.Cells.Item(11, 57).Value = 12
$ws.Cells.Item(11, 58).Value = 3
$ws.Cells.Item(11, 59).Value = 13
$ws.Cells.Item(11, 60).Value = ""

# Row 12: 2020年
$ws.Cells.Item(12, 1).Value = "2020年"
$ws.Cells.Item(12, 2).Value = ""
$ws.Cells.Item(12, 3).Value = 15
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = 6
$ws.Cells.Item(12, 6).Value = 4
$ws.Cells.Item(12, 7).Value = ""
$ws.Cells.Item(12, 8).Value = 14
$ws.Cells.Item(12, 9).Value = ""
$ws.Cells.Item(12, 10).Value = 14
$ws.Cells.Item(12, 11).Value = ""
$ws.Cells.Item(12, 12).Value = 6
$ws.Cells.Item(12, 13).Value = 4
$ws.Cells.Item(12, 14).Value = ""
$ws.Cells.Item(12, 15).Value = 15
$ws.Cells.Item(12, 16).Value = 16
$ws.Cells.Item(12, 17).Value = 17
$ws.Cells.Item(12, 18).Value = 1
$ws.Cells.Item(12, 19).Value = 3
$ws.Cells.Item(12, 20).Value = 5
$ws.Cells.Item(12, 21).Value = 9
$ws.Cells.Item(12, 22).Value = 7
$ws.Cells.Item(12, 23).Value = 8
$ws.Cells.Item(12, 24).Value = 11
$ws.Cells.Item(12, 25).Value = 12
$ws.Cells.Item(12, 26).Value = 2
$ws.Cells.Item(12, 27).Value = 13
$ws.Cells.Item(12, 28).Value = ""
$ws.Cells.Item(12, 29).Value = 17
$ws.Cells.Item(12, 30).Value = 16
$ws.Cells.Item(12, 31).Value = 2
$ws.Cells.Item(12, 32).Value = 1
$ws.Cells.Item(12, 33).Value = 5
$ws.Cells.Item(12, 34).Value = 9
$ws.Cells.Item(12, 35).Value = 7
$ws.Cells.Item(12, 36).Value = ""
$ws.Cells.Item(12, 37).Value = 15
$ws.Cells.Item(12, 38).Value = ""
$ws.Cells.Item(12, 39).Value = 6
$ws.Cells.Item(12, 40).Value = 4
$ws.Cells.Item(12, 41).Value = ""
$ws.Cells.Item(12, 42).Value = 14
$ws.Cells.Item(12, 43).Value = 16
$ws.Cells.Item(12, 44).Value = 2
$ws.Cells.Item(12, 45).Value = 1
$ws.Cells.Item(12, 46).Value = 5
$ws.Cells.Item(12, 47).Value = 9
$ws.Cells.Item(12, 48).Value = 7
$ws.Cells.Item(12, 49).Value = 8
$ws.Cells.Item(12, 50).Value = 11
$ws.Cells.Item(12, 51).Value = 12
$ws.Cells.Item(12, 52).Value = 3
$ws.Cells.Item(12, 53).Value = 13
$ws.Cells.Item(12, 54).Value = ""
$ws.Cells.Item(12, 55).Value = 8
$ws.Cells.Item(12, 56).Value = 11
$ws.Cells.Item(12, 57).Value = 12
$ws.Cells.Item(12, 58).Value = 3
$ws.Cells.Item(12, 59).Value = 13
$ws.Cells.Item(12, 60).Value = ""

# Row 13: 2021年
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 10
$ws.Cells.Item(13, 3).Value = 14
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = 6
$ws.Cells.Item(13, 6).Value = 4
$ws.Cells.Item(13, 7).Value = ""
$ws.Cells.Item(13, 8).Value = 15
$ws.Cells.Item(13, 9).Value = 11
$ws.Cells.Item(13, 10).Value = 14
$ws.Cells.Item(13, 11).Value = ""
$ws.Cells.Item(13, 12).Value = 6
$ws.Cells.Item(13, 13).Value = 4
$ws.Cells.Item(13, 14).Value = ""
$ws.Cells.Item(13, 15).Value = 15
$ws.Cells.Item(13, 16).Value = 17
$ws.Cells.Item(13, 17).Value = 16
$ws.Cells.Item(13, 18).Value = 1
$ws.Cells.Item(13, 19).Value = 3
$ws.Cells.Item(13, 20).Value = 5
$ws.Cells.Item(13, 21).Value = 9
$ws.Cells.Item(13, 22).Value = 7
$ws.Cells.Item(13, 23).Value = 8
$ws.Cells.Item(13, 24).Value = 10
$ws.Cells.Item(13, 25).Value = 12
$ws.Cells.Item(13, 26).Value = 2
$ws.Cells.Item(13, 27).Value = 13
$ws.Cells.Item(13, 28).Value = ""
$ws.Cells.Item(13, 29).Value = 17
$ws.Cells.Item(13, 30).Value = 16
$ws.Cells.Item(13, 31).Value = 1
$ws.Cells.Item(13, 32).Value = 2
$ws.Cells.Item(13, 33).Value = 5
$ws.Cells.Item(13, 34).Value = 9
$ws.Cells.Item(13, 35).Value = 7
$ws.Cells.Item(13, 36).Value = 10
$ws.Cells.Item(13, 37).Value = 14
$ws.Cells.Item(13, 38).Value = ""
$ws.Cells.Item(13, 39).Value = 6
$ws.Cells.Item(13, 40).Value = 4
$ws.Cells.Item(13, 41).Value = ""
$ws.Cells.Item(13, 42).Value = 15
$ws.Cells.Item(13, 43).Value = 16
$ws.Cells.Item(13, 44).Value = 2
$ws.Cells.Item(13, 45).Value = 1
$ws.Cells.Item(13, 46).Value = 5
$ws.Cells.Item(13, 47).Value = 9
$ws.Cells.Item(13, 48).Value = 7
$ws.Cells.Item(13, 49).Value = 8
$ws.Cells.Item(13, 50).Value = 11
$ws.Cells.Item(13, 51).Value = 12
$ws.Cells.Item(13, 52).Value = 3
$ws.Cells.Item(13, 53).Value = 13
$ws.Cells.Item(13, 54).Value = ""
$ws.Cells.Item(13, 55).Value = 8
$ws.Cells.Item(13, 56).Value = 11
$ws.Cells.Item(13, 57).Value = 12
$ws.Cells.Item(13, 58).Value = 3
$ws.Cells.Item(13, 59).Value = 13
$ws.Cells.Item(13, 60).Value = ""
